# Apply updated x_m (E) and y_m (F) sensor coordinate values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; E = 10.48128064354186; F = 2.323561441160277 },
    @{ Row = 3; E = 10.48128064354186; F = 2.323561441160277 },
    @{ Row = 4; E = 10.48128064354186; F = 2.323561441160277 },
    @{ Row = 5; E = 10.48128064354186; F = 2.323561441160277 },
    @{ Row = 6; E = 12.44709339821764; F = 0.8011205862879082 },
    @{ Row = 7; E = 12.44709339821764; F = 0.8011205862879082 },
    @{ Row = 8; E = 12.44709339821764; F = 0.8011205862879082 },
    @{ Row = 9; E = 12.44709339821764; F = 0.8011205862879082 },
    @{ Row = 10; E = 11.46418702087975; F = 1.562341013724093 },
    @{ Row = 11; E = 11.46418702087975; F = 1.562341013724093 },
    @{ Row = 12; E = 11.46418702087975; F = 1.562341013724093 },
    @{ Row = 13; E = 11.46418702087975; F = 1.562341013724093 },
    @{ Row = 14; E = 11.46418702087975; F = 1.562341013724093 },
    @{ Row = 15; E = 11.46418702087975; F = 1.562341013724093 },
    @{ Row = 16; E = 14.41290615289341; F = 2.323561441160277 },
    @{ Row = 17; E = 14.41290615289341; F = 2.323561441160277 },
    @{ Row = 18; E = 14.41290615289341; F = 2.323561441160277 },
    @{ Row = 19; E = 14.41290615289341; F = 2.323561441160277 },
    @{ Row = 20; E = 16.37871890756918; F = 0.8011205862879082 },
    @{ Row = 21; E = 16.37871890756918; F = 0.8011205862879082 },
    @{ Row = 22; E = 16.37871890756918; F = 0.8011205862879082 },
    @{ Row = 23; E = 16.37871890756918; F = 0.8011205862879082 },
    @{ Row = 24; E = 15.39581253023129; F = 1.562341013724093 },
    @{ Row = 25; E = 15.39581253023129; F = 1.562341013724093 },
    @{ Row = 26; E = 15.39581253023129; F = 1.562341013724093 },
    @{ Row = 27; E = 15.39581253023129; F = 1.562341013724093 },
    @{ Row = 28; E = 15.39581253023129; F = 1.562341013724093 },
    @{ Row = 29; E = 15.39581253023129; F = 1.562341013724093 },
    @{ Row = 30; E = 18.34453166224496; F = 2.323561441160277 },
    @{ Row = 31; E = 18.34453166224496; F = 2.323561441160277 },
    @{ Row = 32; E = 18.34453166224496; F = 2.323561441160277 },
    @{ Row = 33; E = 18.34453166224496; F = 2.323561441160277 },
    @{ Row = 34; E = 20.31034441692073; F = 0.8011205862879082 },
    @{ Row = 35; E = 20.31034441692073; F = 0.8011205862879082 },
    @{ Row = 36; E = 20.31034441692073; F = 0.8011205862879082 },
    @{ Row = 37; E = 20.31034441692073; F = 0.8011205862879082 },
    @{ Row = 38; E = 19.32743803958284; F = 1.562341013724093 },
    @{ Row = 39; E = 19.32743803958284; F = 1.562341013724093 },
    @{ Row = 40; E = 19.32743803958284; F = 1.562341013724093 },
    @{ Row = 41; E = 19.32743803958284; F = 1.562341013724093 },
    @{ Row = 42; E = 19.32743803958284; F = 1.562341013724093 },
    @{ Row = 43; E = 19.32743803958284; F = 1.562341013724093 },
    @{ Row = 44; E = 22.88728249497977; F = 1.532865241589286 },
    @{ Row = 45; E = 22.88728249497977; F = 1.532865241589286 },
    @{ Row = 46; E = 3.524707142899675; F = 2.306786678521763 },
    @{ Row = 47; E = 3.524707142899675; F = 2.306786678521763 },
    @{ Row = 48; E = 3.524707142899675; F = 2.306786678521763 },
    @{ Row = 49; E = 3.524707142899675; F = 2.306786678521763 },
    @{ Row = 50; E = 5.313474859126736; F = 0.848166496176657 },
    @{ Row = 51; E = 5.313474859126736; F = 0.848166496176657 },
    @{ Row = 52; E = 5.313474859126736; F = 0.848166496176657 },
    @{ Row = 53; E = 5.313474859126736; F = 0.848166496176657 },
    @{ Row = 54; E = 4.419091001013205; F = 1.57747658734921 },
    @{ Row = 55; E = 4.419091001013205; F = 1.57747658734921 },
    @{ Row = 56; E = 4.419091001013205; F = 1.57747658734921 },
    @{ Row = 57; E = 4.419091001013205; F = 1.57747658734921 },
    @{ Row = 58; E = 4.419091001013205; F = 1.57747658734921 },
    @{ Row = 59; E = 4.419091001013205; F = 1.57747658734921 },
    @{ Row = 60; E = 7.859727870232616; F = 1.595103666437947 },
    @{ Row = 61; E = 7.859727870232616; F = 1.595103666437947 },
    @{ Row = 62; E = 8.30672713918489; F = 2.29948305054856 },
    @{ Row = 63; E = 8.30672713918489; F = 2.29948305054856 },
    @{ Row = 64; E = 7.511669452373797; F = 2.236927454735504 },
    @{ Row = 65; E = 7.511669452373797; F = 2.236927454735504 },
    @{ Row = 66; E = 8.300131735166509; F = 2.333117601348132 },
    @{ Row = 67; E = 8.300131735166509; F = 2.333117601348132 },
    @{ Row = 68; E = 7.862060146081291; F = 0.8491103438901633 },
    @{ Row = 69; E = 7.862060146081291; F = 0.8491103438901633 },
    @{ Row = 70; E = 7.771106541044939; F = 1.278022236493023 },
    @{ Row = 71; E = 7.771106541044939; F = 1.278022236493023 },
    @{ Row = 72; E = 8.070627752880217; F = 0.5250975302862629 },
    @{ Row = 73; E = 8.070627752880217; F = 0.5250975302862629 },
    @{ Row = 74; E = 1.346770009225126; F = 1.545855209035676 },
    @{ Row = 75; E = 7.108791215334779; F = 2.168836182098675 },
    @{ Row = 76; E = 7.108791215334779; F = 2.168836182098675 },
    @{ Row = 77; E = 4.058125708832302; F = 2.2797746497292 },
    @{ Row = 78; E = 4.058125708832302; F = 2.2797746497292 },
    @{ Row = 79; E = 4.058125708832302; F = 2.2797746497292 },
    @{ Row = 80; E = 4.058125708832302; F = 2.2797746497292 },
    @{ Row = 81; E = 8.091803230443098; F = 0.5178044016768293 },
    @{ Row = 82; E = 8.091803230443098; F = 0.5178044016768293 },
    @{ Row = 83; E = 7.044702583584955; F = 1.003754635465525 },
    @{ Row = 84; E = 7.044702583584955; F = 1.003754635465525 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 6).Value = $u.F
}
